$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date in column C for rows 2-13
# from serial 45233 (2023-11-03) to serial 45243 (2023-11-13)
for ($row = 2; $row -le 13; $row++) {
    $ws.Cells.Item($row, 3).Value = 45243
}
